$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08830399999999999
$ws.Range("H2").Value = 0.264912
$ws.Range("I2").Value = 0.04372337970871547
$ws.Range("J2").Value = 0.04372337970871546
$ws.Range("M2").Value = 20.854426
$ws.Range("N2").Value = 62.563278
$ws.Range("O2").Value = 0.1507164072139519
$ws.Range("P2").Value = 0.1507164072139519
$ws.Range("Q2").Value = 1.841529233504
$ws.Range("R2").Value = 16.573763101536
$ws.Range("S2").Value = 0.006589830700949003
$ws.Range("T2").Value = 0.006589830700949002
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08830399999999999
$ws.Range("H3").Value = 0.264912
$ws.Range("I3").Value = 0.04372337970871547
$ws.Range("J3").Value = 0.04372337970871546
$ws.Range("O3").Value = 0.6862909728343718
$ws.Range("P3").Value = 0.6862909728343718
$ws.Range("Q3").Value = 8.385449948858666
$ws.Range("R3").Value = 75.46904953972799
$ws.Range("S3").Value = 0.03000696079590097
$ws.Range("T3").Value = 0.03000696079590097
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08830399999999999
$ws.Range("H4").Value = 0.264912
$ws.Range("I4").Value = 0.04372337970871547
$ws.Range("J4").Value = 0.04372337970871546
$ws.Range("N4").Value = 67.65920700000001
$ws.Range("O4").Value = 0.1629926199516763
$ws.Range("P4").Value = 0.1629926199516763
$ws.Range("Q4").Value = 1.991526204976
$ws.Range("R4").Value = 17.923735844784
$ws.Range("S4").Value = 0.007126588211865494
$ws.Range("T4").Value = 0.007126588211865493
$ws.Range("I5").Value = 0.5310748730197871
$ws.Range("J5").Value = 0.531074873019787
$ws.Range("M5").Value = 20.854426
$ws.Range("N5").Value = 62.563278
$ws.Range("O5").Value = 0.1507164072139519
$ws.Range("P5").Value = 0.1507164072139519
$ws.Range("Q5").Value = 22.367664859412
$ws.Range("R5").Value = 201.308983734708
$ws.Range("S5").Value = 0.08004169682314805
$ws.Range("T5").Value = 0.08004169682314803
$ws.Range("I6").Value = 0.5310748730197871
$ws.Range("J6").Value = 0.531074873019787
$ws.Range("O6").Value = 0.6862909728343718
$ws.Range("P6").Value = 0.6862909728343718
$ws.Range("S6").Value = 0.3644718912526401
$ws.Range("T6").Value = 0.3644718912526401
$ws.Range("I7").Value = 0.5310748730197871
$ws.Range("J7").Value = 0.531074873019787
$ws.Range("N7").Value = 67.65920700000001
$ws.Range("O7").Value = 0.1629926199516763
$ws.Range("P7").Value = 0.1629926199516763
$ws.Range("Q7").Value = 24.189564792778
$ws.Range("S7").Value = 0.08656128494399888
$ws.Range("T7").Value = 0.08656128494399887
$ws.Range("I8").Value = 0.4252017472714976
$ws.Range("J8").Value = 0.4252017472714976
$ws.Range("M8").Value = 20.854426
$ws.Range("N8").Value = 62.563278
$ws.Range("O8").Value = 0.1507164072139519
$ws.Range("P8").Value = 0.1507164072139519
$ws.Range("Q8").Value = 17.90852978324
$ws.Range("R8").Value = 161.17676804916
$ws.Range("S8").Value = 0.06408487968985491
$ws.Range("T8").Value = 0.06408487968985491
$ws.Range("I9").Value = 0.4252017472714976
$ws.Range("J9").Value = 0.4252017472714976
$ws.Range("O9").Value = 0.6862909728343718
$ws.Range("P9").Value = 0.6862909728343718
$ws.Range("S9").Value = 0.2918121207858308
$ws.Range("T9").Value = 0.2918121207858307
$ws.Range("I10").Value = 0.4252017472714976
$ws.Range("J10").Value = 0.4252017472714976
$ws.Range("N10").Value = 67.65920700000001
$ws.Range("O10").Value = 0.1629926199516763
$ws.Range("P10").Value = 0.1629926199516763
$ws.Range("S10").Value = 0.06930474679581192
$ws.Range("T10").Value = 0.0693047467958119